$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 31   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/21/2024  Through  10/27/2024"

# --- Helper: set a cell to a text "N/A" placeholder while keeping the General-format/style
#     used by the other text placeholder cells in this table (copied from A14). ---
function Set-TextPlaceholder($ref, $text) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Range("A14").Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# Row 14
Set-TextPlaceholder "G14" "0"
Set-TextPlaceholder "H14" "***.*"
$ws.Range("M14").Value = -33.333333333333

# Row 15
Set-TextPlaceholder "C15" "0"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = -17.647058823529
$ws.Range("M15").Value = 100

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -31.818181818181
$ws.Range("I16").Value = 255
$ws.Range("J16").Value = 263
$ws.Range("K16").Value = -3.041825095057
$ws.Range("L16").Value = 2.409638554216
$ws.Range("M16").Value = 53.614457831325
$ws.Range("N16").Value = -74.216380182002

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = -18.75
$ws.Range("I17").Value = 425
$ws.Range("J17").Value = 370
$ws.Range("K17").Value = 14.864864864864
$ws.Range("L17").Value = 53.985507246376
$ws.Range("M17").Value = 133.516483516484
$ws.Range("N17").Value = 40.264026402640

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 14
$ws.Range("E18").Value = -78.571428571428
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 48
$ws.Range("H18").Value = -43.75
$ws.Range("I18").Value = 334
$ws.Range("J18").Value = 463
$ws.Range("K18").Value = -27.861771058315
$ws.Range("L18").Value = -26.754385964912
$ws.Range("M18").Value = -17.326732673267
$ws.Range("N18").Value = -83.283283283283

# Row 19
$ws.Range("C19").Value = 25
$ws.Range("D19").Value = 30
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 90
$ws.Range("G19").Value = 107
$ws.Range("H19").Value = -15.887850467289
$ws.Range("I19").Value = 1082
$ws.Range("J19").Value = 1139
$ws.Range("K19").Value = -5.004389815627
$ws.Range("L19").Value = -12.952534191472
$ws.Range("M19").Value = 96.727272727272
$ws.Range("N19").Value = -2.169981916817

# Row 20
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 15
$ws.Range("E20").Value = -46.666666666666
$ws.Range("F20").Value = 46
$ws.Range("G20").Value = 63
$ws.Range("H20").Value = -26.984126984127
$ws.Range("I20").Value = 457
$ws.Range("J20").Value = 431
$ws.Range("K20").Value = 6.032482598607
$ws.Range("L20").Value = 110.599078341014
$ws.Range("M20").Value = 90.416666666666
$ws.Range("N20").Value = -86.40285629277

# Row 21
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 73
$ws.Range("E21").Value = -38.356164383561
$ws.Range("F21").Value = 205
$ws.Range("G21").Value = 273
$ws.Range("H21").Value = -24.908424908424
$ws.Range("I21").Value = 2585
$ws.Range("J21").Value = 2701
$ws.Range("K21").Value = -4.294705664568
$ws.Range("L21").Value = 4.275917708753
$ws.Range("M21").Value = 65.492957746478
$ws.Range("N21").Value = -66.824948665297

# Row 22
$ws.Range("D22").Value = 2
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -83.333333333333
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = -23.333333333333
$ws.Range("L22").Value = -32.352941176470
$ws.Range("M22").Value = 475

# Row 23
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = 26.666666666666

# Row 24
$ws.Range("C24").Value = 53
$ws.Range("D24").Value = 47
$ws.Range("E24").Value = 12.765957446808
$ws.Range("F24").Value = 233
$ws.Range("G24").Value = 212
$ws.Range("H24").Value = 9.905660377358
$ws.Range("I24").Value = 2278
$ws.Range("J24").Value = 2345
$ws.Range("K24").Value = -2.857142857142
$ws.Range("L24").Value = -0.697471665213
$ws.Range("M24").Value = 75.230769230769

# Row 25
$ws.Range("C25").Value = 33
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = 32
$ws.Range("F25").Value = 142
$ws.Range("G25").Value = 115
$ws.Range("H25").Value = 23.478260869565
$ws.Range("I25").Value = 1502
$ws.Range("J25").Value = 1285
$ws.Range("K25").Value = 16.887159533073
$ws.Range("L25").Value = 23.41824157765

# Row 26
$ws.Range("C26").Value = 22
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 37.5
$ws.Range("F26").Value = 88
$ws.Range("G26").Value = 75
$ws.Range("H26").Value = 17.333333333333
$ws.Range("I26").Value = 808
$ws.Range("J26").Value = 740
$ws.Range("K26").Value = 9.189189189189
$ws.Range("L26").Value = 39.550949913644
$ws.Range("M26").Value = 35.798319327731

# Row 27
Set-TextPlaceholder "C27" "0"
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = -20

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -66.666666666666
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 79
$ws.Range("J28").Value = 87
$ws.Range("K28").Value = -9.195402298850
$ws.Range("L28").Value = 17.910447761194

# Row 29
Set-TextPlaceholder "D29" "0"
Set-TextPlaceholder "E29" "***.*"

# Row 30
Set-TextPlaceholder "D30" "0"
Set-TextPlaceholder "E30" "***.*"
